# LOM3257.xlsx update (site rebuild 2023-01-09)
#
# Net effect on the "Ficha" worksheet:
#   - Ativacao date bumped from 01/01/2020 to 01/01/2023
#   - a new "Docentes responsaveis:" row is inserted right after "Objectives:"
#   - the short/long syllabus bodies (English course description) are filled in
#   - the professor responsible for "Metodo" changes
#   - the "Norma de recuperacao" text is replaced with a new grading policy

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write literal text into a cell without letting Excel's "smart"
# input parser reinterpret date-looking text (e.g. "01/01/2023") as a real
# date serial number, and without leaving the cell's own style (s="2"/"s=3")
# altered. We flip a scratch cell to Text format once, then for each target
# cell: stash its current format on a second scratch cell, borrow the Text
# format just long enough to type the value, then paste the original format
# back.
# ---------------------------------------------------------------------------
$fmtScratch = $ws.Range("ZZ1")
$fmtScratch.NumberFormat = "@"
$holdScratch = $ws.Range("ZZ2")

function Set-LiteralText($cell, [string]$text) {
    $cell.Copy() | Out-Null
    $holdScratch.PasteSpecial(-4122) | Out-Null   # xlPasteFormats: stash original format
    $fmtScratch.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null          # xlPasteFormats: borrow Text format
    $cell.Value = $text
    $holdScratch.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null          # xlPasteFormats: restore original format
}

# 1. Ativacao: 01/01/2020 -> 01/01/2023
Set-LiteralText $ws.Range("B8") "01/01/2023"
Set-LiteralText $ws.Range("C8") "01/01/2023"

# 2. Insert a new blank row at 12 ("Docentes responsaveis:"), pushing
#    everything from the old row 12 onward down by one.
$ws.Rows.Item(12).Insert()

# 3. New row 12: label only
$ws.Range("A12").Value = "Docentes responsáveis:"

# 4. Old row 12 ("Programa resumido:" / professor name) is now row 13; its
#    B/C value becomes the (re-used) activation date string.
Set-LiteralText $ws.Range("B13") "01/01/2023"
Set-LiteralText $ws.Range("C13") "01/01/2023"

# 5. Old row 13 ("Short syllabus:") is now row 14; give it a body text.
$shortSyllabus = "Review of Newtonian dynamics (" + [char]0x201C + "vector mechanics" + [char]0x201D + "). General structure of classical mechanics in the Lagrangian and Hamiltonian formulations. Applications to problems of central forces and dynamics of rigid bodies. Problems in non-inertial frames."
$ws.Range("B14").Value = $shortSyllabus
$ws.Range("C14").Value = $shortSyllabus

# 6. Old row 14 ("Programa:") is now row 15; its value becomes the other professor's name.
$ws.Range("B15").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("C15").Value = "5840730 - Antonio Jefferson da Silva Machado"

# 7. Old row 15 ("Syllabus:") is now row 16; give it a body text.
$syllabus = "1. Review of point kinematics: position, velocity and acceleration vectors. 2. Forces, resultant force, conservation of linear momentum and Newton's Laws of Dynamics; 3. Work and energy; kinetic energy and potential energy. Kinetic energy theorem and conservation of total energy. 4. Hamilton's variational principle and the Lagrangian Formulation of Mechanics: generalized coordinates, Lagrangian and the Euler-Lagrange equations. Ignorable coordinates and conservation principles. 5. Applications: central forces; torque and conservation of angular momentum; two-body problems with mutual attraction or repulsion; Gravitation and Kepler's laws. 6. Linearly accelerated and rotating non-inertial frames of reference. Inertia forces: centrifugal force, Coriolis force. Effects of inertia forces on planet Earth; Foucault pendulum. 7. Statics and Dynamics of Rigid Bodies; moments of inertia; torques; Plane motion of rigid bodies; parallel axis theorem. 8. Hamilton's Formulation for Classical Mechanics: the Hamiltonian and Hamilton's equations"
$ws.Range("B16").Value = $syllabus
$ws.Range("C16").Value = $syllabus

# 8. Old row 17 ("Metodo:") is now row 18; professor in charge changed.
$ws.Range("B18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"

# 9. Old row 19 ("Norma de recuperacao:") is now row 20; replace its text.
$recoveryNorm = "A nota final (NF) será a média ponderada de três provas, P1 (peso 1), P2 (peso 1) e P3 (peso 2)"
$ws.Range("B20").Value = $recoveryNorm
$ws.Range("C20").Value = $recoveryNorm

# cleanup scratch cells used for the literal-text trick
$ws.Range("ZZ1:ZZ2").Clear() | Out-Null
